$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("diagnostics")
if (-not $ws) { $ws = $wb.ActiveSheet }

$ws.Range("A2").Value = "00b35d99-0357-428a-b000-b72facf22db3"
$ws.Range("B2").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C2").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A3").Value = "f5354f01-5056-499b-9650-ad0e34c26863"
$ws.Range("B3").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C3").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A4").Value = "e75bfd1c-54cc-4ee4-a89c-66d96d99721e"
$ws.Range("B4").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C4").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A5").Value = "731be840-49de-42d4-8565-381b050d6cdf"
$ws.Range("B5").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C5").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A6").Value = "5f71dfd2-e6e9-4584-ad56-cf1ef5a6d362"
$ws.Range("B6").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C6").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A7").Value = "03780194-cdff-40ca-bd7f-bfe75a98a13a"
$ws.Range("B7").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C7").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A8").Value = "c088fb5f-c642-47ee-a94a-923160c9b8e0"
$ws.Range("B8").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C8").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A9").Value = "c469e10f-2238-4d5e-8331-eecdd585444d"
$ws.Range("B9").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C9").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A10").Value = "6ac837ca-e253-471b-a5c7-6ae567ca5f35"
$ws.Range("B10").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C10").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A11").Value = "ca168835-8342-42de-949e-1784c845e974"
$ws.Range("B11").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C11").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A12").Value = "5bdd1fba-43eb-40f0-9092-fad37fbf389d"
$ws.Range("B12").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C12").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A13").Value = "29fa6165-e7b9-48b3-8d8c-8e6a856212f2"
$ws.Range("B13").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C13").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A14").Value = "d8c262e0-63f5-4259-9b81-241941064f59"
$ws.Range("B14").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C14").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A15").Value = "12fd13c1-5de0-4a55-8147-981ac46c4847"
$ws.Range("B15").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C15").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A16").Value = "cb71e860-f8ce-40bf-8d91-e04bec3d031c"
$ws.Range("B16").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C16").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A17").Value = "7c832afd-6bef-4368-89f5-e298e3805273"
$ws.Range("B17").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C17").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A18").Value = "57d4c2bb-6a85-45e9-9fc8-f7c6e9aa96ec"
$ws.Range("B18").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C18").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A19").Value = "3dc1ccee-0f2b-4cd9-853d-090170a60e9c"
$ws.Range("B19").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C19").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A20").Value = "d71957d5-8dc8-4b91-93c7-869cdd440180"
$ws.Range("B20").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C20").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A21").Value = "e8423fd1-6505-4199-85cf-5bab26b212ba"
$ws.Range("B21").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C21").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A22").Value = "f2f6b6e3-bdcc-4de4-9134-b8565a3c8779"
$ws.Range("B22").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C22").Value = "05e8feaa-0bed-5909-a817-39812494b361"
$ws.Range("A23").Value = "6142938a-80b1-452c-9a21-30771b601f5a"
$ws.Range("B23").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C23").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A24").Value = "f82aff3e-9da3-44a9-a662-ddc0476c662b"
$ws.Range("B24").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C24").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A25").Value = "1d9c3df2-fbfb-475e-bc05-57f3664196c4"
$ws.Range("B25").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C25").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A26").Value = "d7aa7333-816b-496a-bd0c-027a377f3a54"
$ws.Range("B26").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C26").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A27").Value = "1fae3601-41e8-4dbc-a24c-2589d4943893"
$ws.Range("B27").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C27").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A28").Value = "54e063db-ff31-472b-a130-8f943c79d0ad"
$ws.Range("B28").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C28").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A29").Value = "b89cf073-d5a4-4933-8549-6a344a66f9ea"
$ws.Range("B29").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C29").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A30").Value = "3cdd2ab0-38de-4489-84ce-61d8dc293c4b"
$ws.Range("B30").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C30").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A31").Value = "511793bd-7bb1-4000-a755-c9378654d240"
$ws.Range("B31").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C31").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A32").Value = "ffd14f6c-7e41-463d-a212-29ee738156e6"
$ws.Range("B32").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C32").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A33").Value = "5a738940-1c0d-4657-8f88-8872c2d68f9b"
$ws.Range("B33").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C33").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A34").Value = "5f15bb03-fb09-42de-a642-252d37ae0040"
$ws.Range("B34").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C34").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A35").Value = "faf210d8-e8b0-4573-b877-1456861cc9c8"
$ws.Range("B35").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C35").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A36").Value = "3d8878f8-2f86-484d-8bbe-882d2d12af9f"
$ws.Range("B36").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C36").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A37").Value = "288b64e6-3c1c-48af-8917-b1662383cb39"
$ws.Range("B37").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C37").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A38").Value = "d9389775-5d4c-4672-80b7-8b3612d53e22"
$ws.Range("B38").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C38").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A39").Value = "a45be105-d0a6-4c4d-bb8c-d575550f0331"
$ws.Range("B39").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C39").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A40").Value = "80a47034-c25c-4e51-894d-5c8ade23e8d1"
$ws.Range("B40").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C40").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A41").Value = "4c80454c-b849-4f93-8b0c-8728089844e5"
$ws.Range("B41").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C41").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A42").Value = "f7e6d31a-837c-48b7-8927-39e6931da65f"
$ws.Range("B42").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C42").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A43").Value = "2e201f48-6326-4974-bb92-a1344a9e564f"
$ws.Range("B43").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C43").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A44").Value = "da0c3596-3d86-4fbd-9cb4-dc4d93cca4bc"
$ws.Range("B44").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C44").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A45").Value = "b06b0b90-ea1c-47d9-a680-3b954d0b5c09"
$ws.Range("B45").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C45").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A46").Value = "0c298549-5d78-46cd-a0d1-274839e8f474"
$ws.Range("B46").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C46").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A47").Value = "3b36c04b-b1a0-42db-827f-1448ef1fae30"
$ws.Range("B47").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C47").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A48").Value = "690a2a5d-7bd1-44c3-8176-b7d892943f1a"
$ws.Range("B48").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C48").Value = "8f460419-7b80-516d-8919-84520950f612"
$ws.Range("A49").Value = "22002859-7fe2-40ce-94a3-5234c9cdd50c"
$ws.Range("B49").Value = "05269d28-15ae-5bd6-bd88-f949ccfa52d7"
$ws.Range("C49").Value = "8f460419-7b80-516d-8919-84520950f612"
